# Add season-record columns (Wins / Losses / Ties) to the player table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells: copy the formatting used by the rest of the header row
# (bold font, thin border, centered) from AC1, then set the new labels.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row shares the same team season record.
$lastRow = 42
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 30).Value = 88
    $ws.Cells.Item($row, 31).Value = 74
    $ws.Cells.Item($row, 32).Value = 0
}
